# Generate Report for Handback
# Fill in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" columns (I, J, K, P) for the 8763a6bc-a198-467d-8971-bf46ea7a2c6d row (row 7)
# on both the zh-cn and de-de sheets, since a (stale) handback was detected for that file.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de6c6129b4e1dfb1095fd1f174185e8b7cb78aa6/e2e/8763a6bc-a198-467d-8971-bf46ea7a2c6d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/942627fb301999d854dc3b05124d9e2594b5f116/e2e/8763a6bc-a198-467d-8971-bf46ea7a2c6d.md."

# ---------------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# I7: Latest Target File -> hyperlink to the handback markdown file
$wsZh.Cells.Item(7, 9).Value = "8763a6bc-a198-467d-8971-bf46ea7a2c6d.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/942627fb301999d854dc3b05124d9e2594b5f116/e2e/8763a6bc-a198-467d-8971-bf46ea7a2c6d.md", "", "", "8763a6bc-a198-467d-8971-bf46ea7a2c6d.md")
$wsZh.Cells.Item(7, 9).Font.Name = "Calibri"
$wsZh.Cells.Item(7, 9).Font.Underline = 2
$wsZh.Cells.Item(7, 9).Font.Color = 15570276

# J7: Latest Handback File
$wsZh.Cells.Item(7, 10).Value = "8763a6bc-a198-467d-8971-bf46ea7a2c6d.fc41b8f0562f18fe3f6d8aff3718edcd4df3aaa2.zh-cn.xlf"

# K7: Latest Handback DateTime
$wsZh.Cells.Item(7, 11).Value = "2016-08-16 12:54:45"

# P7: Error Detail
$wsZh.Cells.Item(7, 16).Value = $errorMessage

# ---------------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# I7: Latest Target File -> hyperlink to the handback markdown file
$wsDe.Cells.Item(7, 9).Value = "8763a6bc-a198-467d-8971-bf46ea7a2c6d.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/942627fb301999d854dc3b05124d9e2594b5f116/e2e/8763a6bc-a198-467d-8971-bf46ea7a2c6d.md", "", "", "8763a6bc-a198-467d-8971-bf46ea7a2c6d.md")
$wsDe.Cells.Item(7, 9).Font.Name = "Calibri"
$wsDe.Cells.Item(7, 9).Font.Underline = 2
$wsDe.Cells.Item(7, 9).Font.Color = 15570276

# J7: Latest Handback File
$wsDe.Cells.Item(7, 10).Value = "8763a6bc-a198-467d-8971-bf46ea7a2c6d.fc41b8f0562f18fe3f6d8aff3718edcd4df3aaa2.de-de.xlf"

# K7: Latest Handback DateTime
$wsDe.Cells.Item(7, 11).Value = "2016-08-16 12:54:52"

# P7: Error Detail
$wsDe.Cells.Item(7, 16).Value = $errorMessage
